$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume figures per the Jan 3 2023 symbol-list
# refresh. Values are written as literal text (NumberFormat "@" forces
# text storage) to preserve the original "245.28" / "-0.59%" style
# formatting exactly as scraped, rather than letting Excel coerce them
# into numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.59%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.40'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.90%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.252'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.66%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05703'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.44%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.617'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.203'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.32%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8503'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.77%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8939'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.29%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01004'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1,570.61%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.18%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07077'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.10%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03157'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '7.72%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09198'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.92%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001526'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.19%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005925'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.23%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.493'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.09%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.175'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-4.56%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3170'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.36%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03285'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.25%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1276'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.10%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.494'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.53%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04078'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.38%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.06%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.18%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004153'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-17.08%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.82%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.58%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1066'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.52%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003744'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-35.09%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002199'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-9.38%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009147'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '7.73%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005265'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.35%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.01%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1050'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '62.30%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002269'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-10.43%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.01%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.01%'
